# COBRA_structure_fields.xlsx edit
# - switch active sheet/tab from "Programatic Specification" to "Field Specification"
# - change the "g x n" text (Field Specification!B35) to "n x g "
# - give that cell a new left/top aligned style
# - adjust row height for row 35 and column widths on both sheets (autofit-style shrink)

$wb = $excel.ActiveWorkbook

$wsProg = $wb.Worksheets.Item(1)   # "Programatic Specification"
$wsField = $wb.Worksheets.Item(2)  # "Field Specification"

# --- Text fix: "g x n" -> "n x g " ---
$cell = $wsField.Range("B35")
$cell.Value = "n x g "

# --- New style for the fixed cell: left/top alignment (new cellXfs entry) ---
$cell.HorizontalAlignment = -4131   # xlLeft
$cell.VerticalAlignment = -4160     # xlTop

# --- Row 35 height shrinks slightly on the Field Specification sheet ---
$wsField.Rows.Item(35).RowHeight = 13.8

# --- Column widths: both sheets get a touch narrower (cosmetic re-layout) ---
$wsProg.Columns.Item(1).ColumnWidth = 21.57482993197277
$wsProg.Columns.Item(2).ColumnWidth = 8.886054421768707
$wsProg.Columns.Item(3).ColumnWidth = 6.860544217687077
$wsProg.Columns.Item(4).ColumnWidth = 75.03401360544217
$wsProg.Columns.Item(5).ColumnWidth = 21.712585034013568
$wsProg.Columns.Item(6).ColumnWidth = 15.636054421768668
$wsProg.Columns.Item(7).ColumnWidth = 16.447278911564666
$wsProg.Columns.Item(8).ColumnWidth = 12.396258503401366

$wsField.Columns.Item(1).ColumnWidth = 21.57482993197277
$wsField.Columns.Item(2).ColumnWidth = 11.722789115646266
$wsField.Columns.Item(3).ColumnWidth = 42.906462585033964
$wsField.Columns.Item(4).ColumnWidth = 171.54931972789066

# --- Active sheet/tab + selection: workbook now opens on "Field Specification" ---
$wsProg.Activate()
$wsProg.Range("B29").Select()

$wsField.Activate()
$wsField.Range("B35").Select()
